# Apply the commit: "update input excel, add price infomation for Years 2041-2050"
# plus the various view/selection/column-width touch-ups that went along with it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) CoalPlantData sheet: clear the stray "placeholder" values that used to be
#    in A16:A19 (they keep their style, just lose the text), widen column A
#    (best-fit) and move the selection down to A16:G20.
# ---------------------------------------------------------------------------
$wsCoal = $wb.Worksheets.Item("CoalPlantData")
$wsCoal.Range("A16").Value = ""
$wsCoal.Range("A17").Value = ""
$wsCoal.Range("A18").Value = ""
$wsCoal.Range("A19").Value = ""
$wsCoal.Columns.Item(1).AutoFit()
$wsCoal.Columns.Item(1).ColumnWidth = 40
$wsCoal.Range("A16:G20").Select()

# ---------------------------------------------------------------------------
# 2) FC_PPA sheet: also widen column A (best-fit); selection stays the same
#    (A6:A17), the sheet just stops being the active tab (handled below by
#    activating Price_Distribution last).
# ---------------------------------------------------------------------------
$wsFcPpa = $wb.Worksheets.Item("FC_PPA")
$wsFcPpa.Columns.Item(1).AutoFit()
$wsFcPpa.Columns.Item(1).ColumnWidth = 40

# ---------------------------------------------------------------------------
# 3) Price_Distribution sheet: add the Years 2041-2050 price rows (24-33),
#    copying the formatting that row 23 already has, extending the shared
#    "year" formula down column A, and filling in the same constant price
#    figures used for the rest of the table.
# ---------------------------------------------------------------------------
$wsPrice = $wb.Worksheets.Item("Price_Distribution")

# Bring formatting (styles) for the new rows in line with the existing table
$wsPrice.Range("A23:K23").Copy()
$wsPrice.Range("A24:K33").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Column A: incrementing year formula (shared with the A5:A23 formula group)
$wsPrice.Range("A24").Formula = "=A23+1"
$wsPrice.Range("A25").Formula = "=A24+1"
$wsPrice.Range("A26").Formula = "=A25+1"
$wsPrice.Range("A27").Formula = "=A26+1"
$wsPrice.Range("A28").Formula = "=A27+1"
$wsPrice.Range("A29").Formula = "=A28+1"
$wsPrice.Range("A30").Formula = "=A29+1"
$wsPrice.Range("A31").Formula = "=A30+1"
$wsPrice.Range("A32").Formula = "=A31+1"
$wsPrice.Range("A33").Formula = "=A32+1"

# Columns B:K: same constant price figures repeated down every new row
$priceValues = @(2.1908344190476199, 1.61642955640051, 1.3984974952381, 1.28052865296804, 1.1954705631659099, 1.1345047716895, 1.01229800196367, 0.90847234906138896, 0.82517252409944197, 0.674660690005075)

for ($row = 24; $row -le 33; $row++) {
    for ($col = 2; $col -le 11; $col++) {
        $wsPrice.Cells.Item($row, $col).Value = $priceValues[$col - 2]
    }
}

# Selection used while editing this table
$wsPrice.Range("B22:K33").Select()

# Price_Distribution becomes the active/visible tab
$wsPrice.Activate()
